# Auto-generated Excel COM-interop edit script
# Applies the data refresh described in the commit message:
# 'Atualizacao de bases das ligas, do dia: 08-04-2024 as 21:28'

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# =========================================================
# Step 1: refreshed odds for a handful of already-recorded matches
# (upstream source re-emitted these rows with corrected/updated values)
# =========================================================

# Rows 130 and 133 swap content (both kicked off on the same date)
# -- Row 130 --
$ws.Range("B130").Value = 7483189
$ws.Range("F130").Value = "Independiente del Valle"
$ws.Range("G130").Value = "Orense"
$ws.Range("H130").Value = 2
$ws.Range("I130").Value = 2
$ws.Range("K130").Value = 1.4
$ws.Range("L130").Value = 4.75
$ws.Range("M130").Value = 7
$ws.Range("N130").Value = 1.4
$ws.Range("O130").Value = 4.5
$ws.Range("P130").Value = 8
$ws.Range("Q130").Value = -1.25
$ws.Range("R130").Value = 1.875
$ws.Range("S130").Value = 1.925
$ws.Range("U130").Value = 1.925
$ws.Range("V130").Value = 1.875
$ws.Range("X130").Value = 3.5
$ws.Range("AA130").Value = 0.925
$ws.Range("AB130").Value = 0.925
$ws.Range("AC130").Value = -1

# -- Row 133 --
$ws.Range("B133").Value = 7483281
$ws.Range("F133").Value = "SD Aucas"
$ws.Range("G133").Value = "Delfin SC"
$ws.Range("H133").Value = 0
$ws.Range("I133").Value = 0
$ws.Range("K133").Value = 1.909
$ws.Range("L133").Value = 3.25
$ws.Range("M133").Value = 4.2
$ws.Range("N133").Value = 1.909
$ws.Range("O133").Value = 3.5
$ws.Range("P133").Value = 4
$ws.Range("Q133").Value = -0.5
$ws.Range("R133").Value = 1.9
$ws.Range("S133").Value = 1.9
$ws.Range("U133").Value = 1.8
$ws.Range("V133").Value = 2
$ws.Range("X133").Value = 2.5
$ws.Range("AA133").Value = 0.8999999999999999
$ws.Range("AB133").Value = -1
$ws.Range("AC133").Value = 1

# Rows 135, 136 and 137 rotate content (135<-137, 136<-135, 137<-136)
# -- Row 135 --
$ws.Range("B135").Value = 7482832
$ws.Range("F135").Value = "Barcelona Guayaquil"
$ws.Range("G135").Value = "Guayaquil City"
$ws.Range("H135").Value = 2
$ws.Range("I135").Value = 1
$ws.Range("J135").Value = "H"
$ws.Range("K135").Value = 1.363
$ws.Range("L135").Value = 5
$ws.Range("M135").Value = 7.5
$ws.Range("N135").Value = 1.444
$ws.Range("O135").Value = 4
$ws.Range("P135").Value = 8
$ws.Range("Q135").Value = -1.25
$ws.Range("R135").Value = 2.05
$ws.Range("S135").Value = 1.75
$ws.Range("U135").Value = 1.95
$ws.Range("V135").Value = 1.85
$ws.Range("W135").Value = 0.444
$ws.Range("Y135").Value = -1
$ws.Range("Z135").Value = -0.5
$ws.Range("AA135").Value = 0.375
$ws.Range("AB135").Value = 0.95
$ws.Range("AC135").Value = -1

# -- Row 136 --
$ws.Range("B136").Value = 7483188
$ws.Range("F136").Value = "Gualaceo SC"
$ws.Range("G136").Value = "Emelec"
$ws.Range("H136").Value = 0
$ws.Range("K136").Value = 3.6
$ws.Range("L136").Value = 3.3
$ws.Range("M136").Value = 2.05
$ws.Range("N136").Value = 2.6
$ws.Range("O136").Value = 3.25
$ws.Range("P136").Value = 2.75
$ws.Range("Q136").Value = 0
$ws.Range("R136").Value = 1.8
$ws.Range("S136").Value = 2
$ws.Range("U136").Value = 1.975
$ws.Range("V136").Value = 1.825
$ws.Range("Y136").Value = 1.75
$ws.Range("Z136").Value = -1
$ws.Range("AA136").Value = 1
$ws.Range("AB136").Value = -1
$ws.Range("AC136").Value = 0.825

# -- Row 137 --
$ws.Range("B137").Value = 7482867
$ws.Range("F137").Value = "Cumbaya FC"
$ws.Range("G137").Value = "LDU Quito"
$ws.Range("H137").Value = 1
$ws.Range("I137").Value = 2
$ws.Range("J137").Value = "A"
$ws.Range("K137").Value = 5.25
$ws.Range("L137").Value = 3.75
$ws.Range("M137").Value = 1.65
$ws.Range("N137").Value = 9
$ws.Range("O137").Value = 4.5
$ws.Range("P137").Value = 1.363
$ws.Range("Q137").Value = 1.25
$ws.Range("R137").Value = 1.975
$ws.Range("S137").Value = 1.825
$ws.Range("U137").Value = 1.825
$ws.Range("V137").Value = 1.975
$ws.Range("W137").Value = -1
$ws.Range("Y137").Value = 0.363
$ws.Range("Z137").Value = 0.4875
$ws.Range("AA137").Value = -0.5
$ws.Range("AB137").Value = 0.825

# Rows 143 and 145 swap content
# -- Row 143 --
$ws.Range("B143").Value = 7528858
$ws.Range("F143").Value = "Orense"
$ws.Range("G143").Value = "SD Aucas"
$ws.Range("H143").Value = 1
$ws.Range("J143").Value = "A"
$ws.Range("K143").Value = 2.2
$ws.Range("L143").Value = 3.2
$ws.Range("M143").Value = 3.2
$ws.Range("N143").Value = 1.95
$ws.Range("O143").Value = 3.2
$ws.Range("P143").Value = 3.8
$ws.Range("Q143").Value = -0.5
$ws.Range("R143").Value = 1.95
$ws.Range("S143").Value = 1.85
$ws.Range("U143").Value = 1.85
$ws.Range("V143").Value = 1.95
$ws.Range("X143").Value = -1
$ws.Range("Y143").Value = 2.8
$ws.Range("Z143").Value = -1
$ws.Range("AA143").Value = 0.8500000000000001
$ws.Range("AB143").Value = 0.8500000000000001

# -- Row 145 --
$ws.Range("B145").Value = 7528852
$ws.Range("F145").Value = "Delfin SC"
$ws.Range("G145").Value = "Tecnico Universitario"
$ws.Range("H145").Value = 2
$ws.Range("J145").Value = "D"
$ws.Range("K145").Value = 2.1
$ws.Range("L145").Value = 3.4
$ws.Range("M145").Value = 3.1
$ws.Range("N145").Value = 2.1
$ws.Range("O145").Value = 3.4
$ws.Range("P145").Value = 3.1
$ws.Range("Q145").Value = -0.25
$ws.Range("R145").Value = 1.8
$ws.Range("S145").Value = 2
$ws.Range("U145").Value = 1.9
$ws.Range("V145").Value = 1.9
$ws.Range("X145").Value = 2.4
$ws.Range("Y145").Value = -1
$ws.Range("Z145").Value = -0.5
$ws.Range("AA145").Value = 0.5
$ws.Range("AB145").Value = 0.8999999999999999

# =========================================================
# Step 2: a finished match is inserted before the previous last row,
# the (still upcoming) match that used to be last gets its line
# odds refreshed, and six more upcoming matches are appended
# =========================================================

$ws.Rows.Item(198).Insert()

# -- New row 198 (finished match) --
$ws.Range("A198").Value = 196
$ws.Range("B198").Value = 8040562
$ws.Range("C198").Value = "Ecuador LigaPro Serie A"
$ws.Range("D198").Value = "Ecuador LigaPro Serie A"
$ws.Range("E198").Value = 45389.83333333334
$ws.Range("F198").Value = "Tecnico Universitario"
$ws.Range("G198").Value = "Emelec"
$ws.Range("H198").Value = 1
$ws.Range("I198").Value = 1
$ws.Range("J198").Value = "D"
$ws.Range("K198").Value = 2.7
$ws.Range("L198").Value = 3
$ws.Range("M198").Value = 2.8
$ws.Range("N198").Value = 2.8
$ws.Range("O198").Value = 2.9
$ws.Range("P198").Value = 2.8
$ws.Range("Q198").Value = 0
$ws.Range("R198").Value = 1.95
$ws.Range("S198").Value = 1.85
$ws.Range("T198").Value = 2
$ws.Range("U198").Value = 1.775
$ws.Range("V198").Value = 2.025
$ws.Range("W198").Value = -1
$ws.Range("X198").Value = 1.9
$ws.Range("Y198").Value = -1
$ws.Range("Z198").Value = 0
$ws.Range("AA198").Value = -0
$ws.Range("AB198").Value = 0
$ws.Range("AC198").Value = -0
$ws.Range("A198").Style = $ws.Range("A197").Style
$ws.Range("E198").Style = $ws.Range("E197").Style

# -- Row 199 is the previously-last match (shifted down by the insert); --
# -- only the id and the pre-match line need to be refreshed --
$ws.Range("A199").Value = 197
$ws.Range("N199").Value = 3.2
$ws.Range("O199").Value = 3.1
$ws.Range("P199").Value = 2.2
$ws.Range("Q199").Value = 0.25
$ws.Range("R199").Value = 1.825
$ws.Range("S199").Value = 1.975
$ws.Range("U199").Value = 1.9
$ws.Range("V199").Value = 1.9

# -- Six brand-new upcoming matches appended at the bottom --
# -- Row 200 --
$ws.Range("A200").Value = 198
$ws.Range("B200").Value = 7773501
$ws.Range("C200").Value = "Ecuador LigaPro Serie A"
$ws.Range("D200").Value = "Ecuador LigaPro Serie A"
$ws.Range("E200").Value = 45394.875
$ws.Range("F200").Value = "Deportivo Cuenca"
$ws.Range("G200").Value = "SD Aucas"
$ws.Range("K200").Value = 3
$ws.Range("L200").Value = 3.2
$ws.Range("M200").Value = 2.25
$ws.Range("N200").Value = 3.1
$ws.Range("O200").Value = 3.2
$ws.Range("P200").Value = 2.2
$ws.Range("Q200").Value = 0.25
$ws.Range("R200").Value = 1.8
$ws.Range("S200").Value = 2
$ws.Range("T200").Value = 2.5
$ws.Range("U200").Value = 2
$ws.Range("V200").Value = 1.8
$ws.Range("W200").Value = 0
$ws.Range("X200").Value = 0
$ws.Range("Y200").Value = 0
$ws.Range("Z200").Value = 0
$ws.Range("AA200").Value = 0
$ws.Range("A200").Style = $ws.Range("A197").Style
$ws.Range("E200").Style = $ws.Range("E197").Style

# -- Row 201 --
$ws.Range("A201").Value = 199
$ws.Range("B201").Value = 7773503
$ws.Range("C201").Value = "Ecuador LigaPro Serie A"
$ws.Range("D201").Value = "Ecuador LigaPro Serie A"
$ws.Range("E201").Value = 45395.625
$ws.Range("F201").Value = "Mushuc Runa"
$ws.Range("G201").Value = "Tecnico Universitario"
$ws.Range("K201").Value = 2.1
$ws.Range("L201").Value = 3.1
$ws.Range("M201").Value = 3.4
$ws.Range("N201").Value = 2.1
$ws.Range("O201").Value = 3.1
$ws.Range("P201").Value = 3.4
$ws.Range("Q201").Value = -0.25
$ws.Range("R201").Value = 1.8
$ws.Range("S201").Value = 2
$ws.Range("T201").Value = 2.25
$ws.Range("U201").Value = 1.9
$ws.Range("V201").Value = 1.9
$ws.Range("W201").Value = 0
$ws.Range("X201").Value = 0
$ws.Range("Y201").Value = 0
$ws.Range("Z201").Value = 0
$ws.Range("AA201").Value = 0
$ws.Range("A201").Style = $ws.Range("A197").Style
$ws.Range("E201").Style = $ws.Range("E197").Style

# -- Row 202 --
$ws.Range("A202").Value = 200
$ws.Range("B202").Value = 8069719
$ws.Range("C202").Value = "Ecuador LigaPro Serie A"
$ws.Range("D202").Value = "Ecuador LigaPro Serie A"
$ws.Range("E202").Value = 45395.72916666666
$ws.Range("F202").Value = "Macara"
$ws.Range("G202").Value = "Orense"
$ws.Range("K202").Value = 1.95
$ws.Range("L202").Value = 3.25
$ws.Range("M202").Value = 3.5
$ws.Range("N202").Value = 1.909
$ws.Range("O202").Value = 3.3
$ws.Range("P202").Value = 3.6
$ws.Range("Q202").Value = -0.5
$ws.Range("R202").Value = 1.95
$ws.Range("S202").Value = 1.85
$ws.Range("T202").Value = 2.5
$ws.Range("U202").Value = 2
$ws.Range("V202").Value = 1.8
$ws.Range("W202").Value = 0
$ws.Range("X202").Value = 0
$ws.Range("Y202").Value = 0
$ws.Range("Z202").Value = 0
$ws.Range("AA202").Value = 0
$ws.Range("A202").Style = $ws.Range("A197").Style
$ws.Range("E202").Style = $ws.Range("E197").Style

# -- Row 203 --
$ws.Range("A203").Value = 201
$ws.Range("B203").Value = 8069537
$ws.Range("C203").Value = "Ecuador LigaPro Serie A"
$ws.Range("D203").Value = "Ecuador LigaPro Serie A"
$ws.Range("E203").Value = 45395.83333333334
$ws.Range("F203").Value = "Emelec"
$ws.Range("G203").Value = "Cumbaya FC"
$ws.Range("K203").Value = 1.28
$ws.Range("L203").Value = 5.5
$ws.Range("M203").Value = 8.5
$ws.Range("N203").Value = 1.3
$ws.Range("O203").Value = 5.25
$ws.Range("P203").Value = 8.5
$ws.Range("Q203").Value = -1.5
$ws.Range("R203").Value = 1.95
$ws.Range("S203").Value = 1.85
$ws.Range("T203").Value = 2.75
$ws.Range("U203").Value = 1.85
$ws.Range("V203").Value = 1.95
$ws.Range("W203").Value = 0
$ws.Range("X203").Value = 0
$ws.Range("Y203").Value = 0
$ws.Range("Z203").Value = 0
$ws.Range("AA203").Value = 0
$ws.Range("A203").Style = $ws.Range("A197").Style
$ws.Range("E203").Style = $ws.Range("E197").Style

# -- Row 204 --
$ws.Range("A204").Value = 202
$ws.Range("B204").Value = 7773067
$ws.Range("C204").Value = "Ecuador LigaPro Serie A"
$ws.Range("D204").Value = "Ecuador LigaPro Serie A"
$ws.Range("E204").Value = 45396.625
$ws.Range("F204").Value = "Universidad Catolica del Ecuador"
$ws.Range("G204").Value = "Independiente del Valle"
$ws.Range("K204").Value = 2.3
$ws.Range("L204").Value = 3.1
$ws.Range("M204").Value = 3
$ws.Range("N204").Value = 2.5
$ws.Range("O204").Value = 3.1
$ws.Range("P204").Value = 2.7
$ws.Range("Q204").Value = 0
$ws.Range("R204").Value = 1.825
$ws.Range("S204").Value = 1.975
$ws.Range("T204").Value = 2.25
$ws.Range("U204").Value = 1.9
$ws.Range("V204").Value = 1.9
$ws.Range("W204").Value = 0
$ws.Range("X204").Value = 0
$ws.Range("Y204").Value = 0
$ws.Range("Z204").Value = 0
$ws.Range("AA204").Value = 0
$ws.Range("A204").Style = $ws.Range("A197").Style
$ws.Range("E204").Style = $ws.Range("E197").Style

# -- Row 205 --
$ws.Range("A205").Value = 203
$ws.Range("B205").Value = 8069721
$ws.Range("C205").Value = "Ecuador LigaPro Serie A"
$ws.Range("D205").Value = "Ecuador LigaPro Serie A"
$ws.Range("E205").Value = 45396.83333333334
$ws.Range("F205").Value = "Club Atletico Libertad"
$ws.Range("G205").Value = "Barcelona Guayaquil"
$ws.Range("K205").Value = 4
$ws.Range("L205").Value = 3.4
$ws.Range("M205").Value = 1.85
$ws.Range("N205").Value = 3.8
$ws.Range("O205").Value = 3.4
$ws.Range("P205").Value = 1.85
$ws.Range("Q205").Value = 0.5
$ws.Range("R205").Value = 1.9
$ws.Range("S205").Value = 1.9
$ws.Range("T205").Value = 2.5
$ws.Range("U205").Value = 2
$ws.Range("V205").Value = 1.8
$ws.Range("W205").Value = 0
$ws.Range("X205").Value = 0
$ws.Range("Y205").Value = 0
$ws.Range("Z205").Value = 0
$ws.Range("AA205").Value = 0
$ws.Range("A205").Style = $ws.Range("A197").Style
$ws.Range("E205").Style = $ws.Range("E197").Style

